$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prueba")

# --- Row 2 ---
$ws.Range("B2").Value = 40
$ws.Range("C2").Value = 45942
$ws.Range("E2").Value = "A,B"

# --- Row 3 ---
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 45972
$ws.Range("E3").Value = "A,B"

# --- Row 4 ---
$ws.Range("B4").Value = 14
$ws.Range("C4").Value = 45980
$ws.Range("E4").Value = "C,D"

# --- Rows 5-8: clear the registro/no muestras/date/aplican data (date cell stays, but empty) ---
$ws.Range("A5").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()

$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()

$ws.Range("A7").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("E7").ClearContents()

$ws.Range("A8").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("E8").ClearContents()

# --- Rows 12-14: add an empty styled (date-formatted) D cell ---
$ws.Range("D12").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("D13").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("D14").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# --- Rows 15-18: remove entirely (shrinks used range back to row 14) ---
$ws.Range("C15:D18").Clear()

# --- Selection / active sheet bookkeeping ---
$ws.Activate() | Out-Null
$ws.Range("C5").Select() | Out-Null
